# Simplify database schema migration
# Appends a new row (row 76) of data to the bottom of each of the four
# worksheets (DE_LFT_#1, DE_LFT_#2, DE_PLT_#1, DE_PLT_#2), matching the
# existing table schema (columns A..I).

$wb = $excel.ActiveWorkbook

# Row 76 data for each worksheet, keyed by sheet name.
# Columns: A(date serial) B C D E(hex strings) F G H I (numbers)
$rowsData = @{
    "DE_LFT_#1" = @{
        A = 45862.43560185185
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x40"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 320
        I = 14
    }
    "DE_LFT_#2" = @{
        A = 45862.43560185185
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x40"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 320
        I = 14
    }
    "DE_PLT_#1" = @{
        A = 45862.43560185185
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x79"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 121
        I = 7
    }
    "DE_PLT_#2" = @{
        A = 45862.43560185185
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x79"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 121
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rowsData[$ws.Name]
    $newRow = 76

    # Column A keeps the same date/time number format as the row above it.
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
    $ws.Cells.Item($newRow, 1).Value = $data.A

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
